$d = $word.ActiveDocument

# The target paragraph holds the text " Simulation of SHA-1 on Virtual Lab "
# (wrapped in non-breaking spaces) inside the "Output" table cell. The same
# plain-text phrase also occurs earlier in the document (inside the
# "Problem Definition" cell, as part of a longer sentence), so we search
# for the nbsp-wrapped form to land on the unique, correct paragraph.
$nbsp = [char]0x00A0
$needle = $nbsp + "Simulation of SHA-1 on Virtual Lab" + $nbsp

$searchRange = $d.Content
$found = $searchRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target 'Simulation of SHA-1 on Virtual Lab' paragraph."
}

# searchRange now spans just the matched text (not the trailing paragraph
# mark). Move one character past it to land on the start of the following
# paragraph - i.e. the insertion point right after the target paragraph.
$insertPos = $searchRange.End + 1

# Empty "ListParagraph"-styled paragraph, bold/bCs, sz 24 / szCs 24 - same
# shape as the five new blank paragraphs added by the edit.
$blankParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'

for ($i = 0; $i -lt 5; $i++) {
    $r = $d.Range($insertPos, $insertPos)
    [void]$r.InsertXML($blankParaXml)
}
